# Jenkins Changes in scripts
$wb = $excel.ActiveWorkbook

# EnvDetails: switch the configured run from India (Opco 1707) to China (Opco 1284)
$envDetails = $wb.Worksheets.Item("EnvDetails")
$envDetails.Range("B4").Value = "China"
$envDetails.Range("B6").Value = 1284
$envDetails.Range("B6").Select() | Out-Null

# OpcoMapping: record the new China Opco code (1284) and make this the active tab
$opcoMapping = $wb.Worksheets.Item("OpcoMapping")
$opcoMapping.Range("B3").Value = 1284

# JIRA_Details: add two new Test Summary rows (global block client/brand),
# reusing the same JIRA references as the preceding "BlockUser" row
$jiraDetails = $wb.Worksheets.Item("JIRA_Details")

$jiraDetails.Range("A34").Value = "BlockGlobalClient"
$jiraDetails.Range("B34").Value = "TSTAUTO-58"
$jiraDetails.Range("C34").Value = "TSTAUTO-22"
$jiraDetails.Range("B34").Font.Color = 0
$jiraDetails.Range("C34").Font.Color = 0

$jiraDetails.Range("A35").Value = "BlockGlobalBrand"
$jiraDetails.Range("B35").Value = "TSTAUTO-58"
$jiraDetails.Range("C35").Value = "TSTAUTO-22"
$jiraDetails.Range("B35").Font.Color = 0
$jiraDetails.Range("C35").Font.Color = 0

$jiraDetails.Activate()
$excel.ActiveWindow.ScrollRow = 19
$jiraDetails.Range("B34:C35").Select() | Out-Null

# OpcoMapping ends up as the selected/active tab in the saved workbook
$opcoMapping.Activate()
$opcoMapping.Range("B3").Select() | Out-Null
